$wb = $excel.ActiveWorkbook

# Rename "Students" sheet to "Sheet 1" and remove the "Schools" sheet,
# mirroring the workbook-level sheet-list change in the diff.
$ws = $wb.Worksheets.Item("Students")
$ws.Name = "Sheet 1"

$schools = $wb.Worksheets.Item("Schools")
$schools.Delete()

# Update rows 24-70 of the (now single) data sheet to reflect the new
# "fake scale" rows + re-laid-out pvkat_1..pvkat_5 blocks.
$ws.Cells.Item(24, 1).Value = "skala_fake_item1"
$ws.Cells.Item(24, 2).Value = -97
$ws.Cells.Item(24, 3).Value = "ja"
$ws.Cells.Item(24, 4).Value = ""
$ws.Cells.Item(24, 5).Value = "nein"
$ws.Cells.Item(25, 1).Value = "skala_fake_item1"
$ws.Cells.Item(25, 2).Value = -98
$ws.Cells.Item(25, 3).Value = "ja"
$ws.Cells.Item(25, 4).Value = ""
$ws.Cells.Item(25, 5).Value = "nein"
$ws.Cells.Item(26, 1).Value = "skala_fake_item1"
$ws.Cells.Item(26, 2).Value = -99
$ws.Cells.Item(26, 3).Value = "ja"
$ws.Cells.Item(26, 4).Value = ""
$ws.Cells.Item(26, 5).Value = "nein"
$ws.Cells.Item(27, 1).Value = "skala_fake_item2"
$ws.Cells.Item(27, 2).Value = 1
$ws.Cells.Item(27, 3).Value = "nein"
$ws.Cells.Item(27, 4).Value = "keinerlei Ahnung"
$ws.Cells.Item(27, 5).Value = "nein"
$ws.Cells.Item(28, 1).Value = "skala_fake_item2"
$ws.Cells.Item(28, 2).Value = 2
$ws.Cells.Item(28, 3).Value = "nein"
$ws.Cells.Item(28, 4).Value = "kaum Ahnung"
$ws.Cells.Item(28, 5).Value = "nein"
$ws.Cells.Item(29, 1).Value = "skala_fake_item2"
$ws.Cells.Item(29, 2).Value = 3
$ws.Cells.Item(29, 3).Value = "nein"
$ws.Cells.Item(29, 4).Value = "etwas Ahnung"
$ws.Cells.Item(29, 5).Value = "nein"
$ws.Cells.Item(30, 1).Value = "skala_fake_item2"
$ws.Cells.Item(30, 2).Value = 4
$ws.Cells.Item(30, 3).Value = "nein"
$ws.Cells.Item(30, 4).Value = "viel Ahnung"
$ws.Cells.Item(30, 5).Value = "nein"
$ws.Cells.Item(31, 1).Value = "skala_fake_item2"
$ws.Cells.Item(31, 2).Value = -97
$ws.Cells.Item(31, 3).Value = "ja"
$ws.Cells.Item(31, 4).Value = ""
$ws.Cells.Item(31, 5).Value = "nein"
$ws.Cells.Item(32, 1).Value = "skala_fake_item2"
$ws.Cells.Item(32, 2).Value = -98
$ws.Cells.Item(32, 3).Value = "ja"
$ws.Cells.Item(32, 4).Value = ""
$ws.Cells.Item(32, 5).Value = "nein"
$ws.Cells.Item(33, 1).Value = "skala_fake_item2"
$ws.Cells.Item(33, 2).Value = -99
$ws.Cells.Item(33, 3).Value = "ja"
$ws.Cells.Item(33, 4).Value = ""
$ws.Cells.Item(33, 5).Value = "nein"
$ws.Cells.Item(34, 1).Value = "skala_fake_item3"
$ws.Cells.Item(34, 2).Value = 1
$ws.Cells.Item(34, 3).Value = "nein"
$ws.Cells.Item(34, 4).Value = "keinerlei Ahnung"
$ws.Cells.Item(34, 5).Value = "nein"
$ws.Cells.Item(35, 1).Value = "skala_fake_item3"
$ws.Cells.Item(35, 2).Value = 2
$ws.Cells.Item(35, 3).Value = "nein"
$ws.Cells.Item(35, 4).Value = "kaum Ahnung"
$ws.Cells.Item(35, 5).Value = "nein"
$ws.Cells.Item(36, 1).Value = "skala_fake_item3"
$ws.Cells.Item(36, 2).Value = 3
$ws.Cells.Item(36, 3).Value = "nein"
$ws.Cells.Item(36, 4).Value = "etwas Ahnung"
$ws.Cells.Item(36, 5).Value = "nein"
$ws.Cells.Item(37, 1).Value = "skala_fake_item3"
$ws.Cells.Item(37, 2).Value = 4
$ws.Cells.Item(37, 3).Value = "nein"
$ws.Cells.Item(37, 4).Value = "viel Ahnung"
$ws.Cells.Item(37, 5).Value = "nein"
$ws.Cells.Item(38, 1).Value = "skala_fake_item3"
$ws.Cells.Item(38, 2).Value = -97
$ws.Cells.Item(38, 3).Value = "ja"
$ws.Cells.Item(38, 4).Value = ""
$ws.Cells.Item(38, 5).Value = "nein"
$ws.Cells.Item(39, 1).Value = "skala_fake_item3"
$ws.Cells.Item(39, 2).Value = -98
$ws.Cells.Item(39, 3).Value = "ja"
$ws.Cells.Item(39, 4).Value = ""
$ws.Cells.Item(39, 5).Value = "nein"
$ws.Cells.Item(40, 1).Value = "skala_fake_item3"
$ws.Cells.Item(40, 2).Value = -99
$ws.Cells.Item(40, 3).Value = "ja"
$ws.Cells.Item(40, 4).Value = ""
$ws.Cells.Item(40, 5).Value = "nein"
$ws.Cells.Item(41, 1).Value = "pvkat_1"
$ws.Cells.Item(41, 2).Value = 1
$ws.Cells.Item(41, 3).Value = "nein"
$ws.Cells.Item(41, 4).Value = "Kompetenzstufe 1"
$ws.Cells.Item(41, 5).Value = "nein"
$ws.Cells.Item(42, 1).Value = "pvkat_1"
$ws.Cells.Item(42, 2).Value = 2
$ws.Cells.Item(42, 3).Value = "nein"
$ws.Cells.Item(42, 4).Value = "Kompetenzstufe 2"
$ws.Cells.Item(42, 5).Value = "nein"
$ws.Cells.Item(43, 1).Value = "pvkat_1"
$ws.Cells.Item(43, 2).Value = 3
$ws.Cells.Item(43, 3).Value = "nein"
$ws.Cells.Item(43, 4).Value = "Kompetenzstufe 3"
$ws.Cells.Item(43, 5).Value = "nein"
$ws.Cells.Item(44, 1).Value = "pvkat_1"
$ws.Cells.Item(44, 2).Value = 4
$ws.Cells.Item(44, 3).Value = "nein"
$ws.Cells.Item(44, 4).Value = "Kompetenzstufe 4"
$ws.Cells.Item(44, 5).Value = "nein"
$ws.Cells.Item(45, 1).Value = "pvkat_pooled"
$ws.Cells.Item(45, 2).Value = 1
$ws.Cells.Item(45, 3).Value = "nein"
$ws.Cells.Item(45, 4).Value = "Kompetenzstufe 1"
$ws.Cells.Item(45, 5).Value = "nein"
$ws.Cells.Item(46, 1).Value = "pvkat_pooled"
$ws.Cells.Item(46, 2).Value = 2
$ws.Cells.Item(46, 3).Value = "nein"
$ws.Cells.Item(46, 4).Value = "Kompetenzstufe 2"
$ws.Cells.Item(46, 5).Value = "nein"
$ws.Cells.Item(47, 1).Value = "pvkat_pooled"
$ws.Cells.Item(47, 2).Value = 3
$ws.Cells.Item(47, 3).Value = "nein"
$ws.Cells.Item(47, 4).Value = "Kompetenzstufe 3"
$ws.Cells.Item(47, 5).Value = "nein"
$ws.Cells.Item(48, 1).Value = "pvkat_pooled"
$ws.Cells.Item(48, 2).Value = 4
$ws.Cells.Item(48, 3).Value = "nein"
$ws.Cells.Item(48, 4).Value = "Kompetenzstufe 4"
$ws.Cells.Item(48, 5).Value = "nein"
$ws.Cells.Item(49, 1).Value = "pvkat_pooled"
$ws.Cells.Item(49, 2).Value = 5
$ws.Cells.Item(49, 3).Value = "nein"
$ws.Cells.Item(49, 4).Value = "Kompetenzstufe 5"
$ws.Cells.Item(49, 5).Value = "nein"
$ws.Cells.Item(50, 1).Value = "pvkat_1"
$ws.Cells.Item(50, 2).Value = 5
$ws.Cells.Item(50, 3).Value = "nein"
$ws.Cells.Item(50, 4).Value = "Kompetenzstufe 5"
$ws.Cells.Item(50, 5).Value = "nein"
$ws.Cells.Item(51, 1).Value = "pvkat_2"
$ws.Cells.Item(51, 2).Value = 1
$ws.Cells.Item(51, 3).Value = "nein"
$ws.Cells.Item(51, 4).Value = "Kompetenzstufe 1"
$ws.Cells.Item(51, 5).Value = "nein"
$ws.Cells.Item(52, 1).Value = "pvkat_2"
$ws.Cells.Item(52, 2).Value = 2
$ws.Cells.Item(52, 3).Value = "nein"
$ws.Cells.Item(52, 4).Value = "Kompetenzstufe 2"
$ws.Cells.Item(52, 5).Value = "nein"
$ws.Cells.Item(53, 1).Value = "pvkat_2"
$ws.Cells.Item(53, 2).Value = 3
$ws.Cells.Item(53, 3).Value = "nein"
$ws.Cells.Item(53, 4).Value = "Kompetenzstufe 3"
$ws.Cells.Item(53, 5).Value = "nein"
$ws.Cells.Item(54, 1).Value = "pvkat_2"
$ws.Cells.Item(54, 2).Value = 4
$ws.Cells.Item(54, 3).Value = "nein"
$ws.Cells.Item(54, 4).Value = "Kompetenzstufe 4"
$ws.Cells.Item(54, 5).Value = "nein"
$ws.Cells.Item(55, 1).Value = "pvkat_2"
$ws.Cells.Item(55, 2).Value = 5
$ws.Cells.Item(55, 3).Value = "nein"
$ws.Cells.Item(55, 4).Value = "Kompetenzstufe 5"
$ws.Cells.Item(55, 5).Value = "nein"
$ws.Cells.Item(56, 1).Value = "pvkat_3"
$ws.Cells.Item(56, 2).Value = 1
$ws.Cells.Item(56, 3).Value = "nein"
$ws.Cells.Item(56, 4).Value = "Kompetenzstufe 1"
$ws.Cells.Item(56, 5).Value = "nein"
$ws.Cells.Item(57, 1).Value = "pvkat_3"
$ws.Cells.Item(57, 2).Value = 2
$ws.Cells.Item(57, 3).Value = "nein"
$ws.Cells.Item(57, 4).Value = "Kompetenzstufe 2"
$ws.Cells.Item(57, 5).Value = "nein"
$ws.Cells.Item(58, 1).Value = "pvkat_3"
$ws.Cells.Item(58, 2).Value = 3
$ws.Cells.Item(58, 3).Value = "nein"
$ws.Cells.Item(58, 4).Value = "Kompetenzstufe 3"
$ws.Cells.Item(58, 5).Value = "nein"
$ws.Cells.Item(59, 1).Value = "pvkat_3"
$ws.Cells.Item(59, 2).Value = 4
$ws.Cells.Item(59, 3).Value = "nein"
$ws.Cells.Item(59, 4).Value = "Kompetenzstufe 4"
$ws.Cells.Item(59, 5).Value = "nein"
$ws.Cells.Item(60, 1).Value = "pvkat_3"
$ws.Cells.Item(60, 2).Value = 5
$ws.Cells.Item(60, 3).Value = "nein"
$ws.Cells.Item(60, 4).Value = "Kompetenzstufe 5"
$ws.Cells.Item(60, 5).Value = "nein"
$ws.Cells.Item(61, 1).Value = "pvkat_4"
$ws.Cells.Item(61, 2).Value = 1
$ws.Cells.Item(61, 3).Value = "nein"
$ws.Cells.Item(61, 4).Value = "Kompetenzstufe 1"
$ws.Cells.Item(61, 5).Value = "nein"
$ws.Cells.Item(62, 1).Value = "pvkat_4"
$ws.Cells.Item(62, 2).Value = 2
$ws.Cells.Item(62, 3).Value = "nein"
$ws.Cells.Item(62, 4).Value = "Kompetenzstufe 2"
$ws.Cells.Item(62, 5).Value = "nein"
$ws.Cells.Item(63, 1).Value = "pvkat_4"
$ws.Cells.Item(63, 2).Value = 3
$ws.Cells.Item(63, 3).Value = "nein"
$ws.Cells.Item(63, 4).Value = "Kompetenzstufe 3"
$ws.Cells.Item(63, 5).Value = "nein"
$ws.Cells.Item(64, 1).Value = "pvkat_4"
$ws.Cells.Item(64, 2).Value = 4
$ws.Cells.Item(64, 3).Value = "nein"
$ws.Cells.Item(64, 4).Value = "Kompetenzstufe 4"
$ws.Cells.Item(64, 5).Value = "nein"
$ws.Cells.Item(65, 1).Value = "pvkat_4"
$ws.Cells.Item(65, 2).Value = 5
$ws.Cells.Item(65, 3).Value = "nein"
$ws.Cells.Item(65, 4).Value = "Kompetenzstufe 5"
$ws.Cells.Item(65, 5).Value = "nein"
$ws.Cells.Item(66, 1).Value = "pvkat_5"
$ws.Cells.Item(66, 2).Value = 1
$ws.Cells.Item(66, 3).Value = "nein"
$ws.Cells.Item(66, 4).Value = "Kompetenzstufe 1"
$ws.Cells.Item(66, 5).Value = "nein"
$ws.Cells.Item(67, 1).Value = "pvkat_5"
$ws.Cells.Item(67, 2).Value = 2
$ws.Cells.Item(67, 3).Value = "nein"
$ws.Cells.Item(67, 4).Value = "Kompetenzstufe 2"
$ws.Cells.Item(67, 5).Value = "nein"
$ws.Cells.Item(68, 1).Value = "pvkat_5"
$ws.Cells.Item(68, 2).Value = 3
$ws.Cells.Item(68, 3).Value = "nein"
$ws.Cells.Item(68, 4).Value = "Kompetenzstufe 3"
$ws.Cells.Item(68, 5).Value = "nein"
$ws.Cells.Item(69, 1).Value = "pvkat_5"
$ws.Cells.Item(69, 2).Value = 4
$ws.Cells.Item(69, 3).Value = "nein"
$ws.Cells.Item(69, 4).Value = "Kompetenzstufe 4"
$ws.Cells.Item(69, 5).Value = "nein"
$ws.Cells.Item(70, 1).Value = "pvkat_5"
$ws.Cells.Item(70, 2).Value = 5
$ws.Cells.Item(70, 3).Value = "nein"
$ws.Cells.Item(70, 4).Value = "Kompetenzstufe 5"
$ws.Cells.Item(70, 5).Value = "nein"
